$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Pre-format cells whose new text would otherwise be auto-parsed as a number
# (Excel COM type-inference on Range.Value), so the literal text is preserved exactly.
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D47').NumberFormat = "@"

# Apply the updated cell values (row-by-row, matching the commit diff).
$ws.Range('D2').Value = '33.805.97'
$ws.Range('E2').Value = '  -0.92%  '
$ws.Range('D3').Value = '1.780.14'
$ws.Range('E3').Value = '  -0.97%  '
$ws.Range('E4').Value = '  +0.05%  '
$ws.Range('D5').Value = '224.12'
$ws.Range('E5').Value = '  +0.62%  '
$ws.Range('D6').Value = '0.545'
$ws.Range('E6').Value = '  -1.03%  '
$ws.Range('E7').Value = '  +0.08%  '
$ws.Range('D8').Value = '31.97'
$ws.Range('E8').Value = '  -1.38%  '
$ws.Range('D9').Value = '0.288'
$ws.Range('E9').Value = '  +1.16%  '
$ws.Range('D10').Value = '0.0676'
$ws.Range('E10').Value = '  -5.68%  '
$ws.Range('D11').Value = '0.0935'
$ws.Range('E11').Value = '  +1.26%  '
$ws.Range('D12').Value = '2.035.41'
$ws.Range('E12').Value = '  -0.97%  '
$ws.Range('D13').Value = '11.20'
$ws.Range('E13').Value = '  +4.39%  '
$ws.Range('D14').Value = '1.768.43'
$ws.Range('E14').Value = '  -1.60%  '
$ws.Range('D15').Value = '33.838.45'
$ws.Range('E15').Value = '  -0.85%  '
$ws.Range('E16').Value = '  -3.24%  '
$ws.Range('D17').Value = '4.13'
$ws.Range('E17').Value = '  -2.03%  '
$ws.Range('E18').Value = '  -2.41%  '
$ws.Range('D19').Value = '238.33'
$ws.Range('E19').Value = '  -3.38%  '
$ws.Range('E20').Value = '  -1.78%  '
$ws.Range('E21').Value = '  -0.01%  '
$ws.Range('D22').Value = '10.56'
$ws.Range('E22').Value = '  -2.79%  '
$ws.Range('E23').Value = '  -2.16%  '
$ws.Range('E24').Value = '  -2.66%  '
$ws.Range('E25').Value = '  +0.69%  '
$ws.Range('B26').Value = 'Cosmos'
$ws.Range('C26').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D26').Value = '7.01'
$ws.Range('E26').Value = '  -0.71%  '
$ws.Range('B27').Value = 'EthereumClassic'
$ws.Range('C27').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D27').Value = '16.06'
$ws.Range('E27').Value = '  -2.94%  '
$ws.Range('E28').Value = '  -0.47%  '
$ws.Range('E29').Value = '  +0.12%  '
$ws.Range('E31').Value = '  -2.29%  '
$ws.Range('E32').Value = '  -3.79%  '
$ws.Range('D33').Value = '3.50'
$ws.Range('E33').Value = '  -0.41%  '
$ws.Range('E34').Value = '  -2.19%  '
$ws.Range('D35').Value = '1.382.21'
$ws.Range('E35').Value = '  -2.16%  '
$ws.Range('E36').Value = '  -1.51%  '
$ws.Range('D37').Value = '1.04'
$ws.Range('E37').Value = '  -2.14%  '
$ws.Range('D38').Value = '0.0185'
$ws.Range('E38').Value = '  -1.16%  '
$ws.Range('D39').Value = '2.40'
$ws.Range('E39').Value = '  +2.28%  '
$ws.Range('E40').Value = '  +4.54%  '
$ws.Range('D41').Value = '78.33'
$ws.Range('E41').Value = '  -2.54%  '
$ws.Range('E42').Value = '  -3.95%  '
$ws.Range('E43').Value = '  +11.79%  '
$ws.Range('E44').Value = '  -2.81%  '
$ws.Range('D45').Value = '0.0₆0143'
$ws.Range('E45').Value = '  +16.93%  '
$ws.Range('E46').Value = '  +1.73%  '
$ws.Range('D47').Value = '1.07'
$ws.Range('E47').Value = '  +3.01%  '
$ws.Range('E48').Value = '  +0.73%  '
$ws.Range('E49').Value = '  -2.03%  '
$ws.Range('D50').Value = '1.937.19'
$ws.Range('E50').Value = '  -0.99%  '
